# Add a "Save" column (H) to the s_vals sheet, matching the existing
# header style used by the other header cells (e.g. G1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy G1's formatting (bold/centered/bordered header style) onto H1
# before setting its value, so H1 reuses the same cell style as the
# other header cells instead of creating a brand-new style entry.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New "Save" data column values for the two data rows.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
